$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading-percent results for Case_4_170 (380 kV case)
$ws.Range("B2").Value = 16.1336378423624
$ws.Range("C2").Value = 12.85433496369821
$ws.Range("D2").Value = 14.30960359890199
$ws.Range("E2").Value = 15.15085222293271
$ws.Range("G2").Value = 51.40380372964874
$ws.Range("H2").Value = 19.59928824244815
$ws.Range("I2").Value = 27.75036052041423
$ws.Range("J2").Value = 8.954229939833217
$ws.Range("M2").Value = 20.05577509288377

$ws.Range("B3").Value = 15.72147864068391
$ws.Range("C3").Value = 12.46859528841875
$ws.Range("D3").Value = 14.30078739159286
$ws.Range("E3").Value = 15.17051172976842
$ws.Range("G3").Value = 51.1327564781549
$ws.Range("H3").Value = 19.6158682217605
$ws.Range("I3").Value = 27.79229541746646
$ws.Range("J3").Value = 8.973996717003899
$ws.Range("M3").Value = 19.94146821478699

$ws.Range("B4").Value = 15.4670577010663
$ws.Range("C4").Value = 12.22929965528447
$ws.Range("D4").Value = 14.29834141684143
$ws.Range("E4").Value = 15.18522033258412
$ws.Range("G4").Value = 50.98285864063714
$ws.Range("H4").Value = 19.63052755397703
$ws.Range("I4").Value = 27.82486119810276
$ws.Range("J4").Value = 8.986876989298135
$ws.Range("M4").Value = 19.87549397728861

$ws.Range("B5").Value = 15.36320748332538
$ws.Range("C5").Value = 12.13133521992943
$ws.Range("D5").Value = 14.29809124057
$ws.Range("E5").Value = 15.19187646700408
$ws.Range("G5").Value = 50.92596484706411
$ws.Range("H5").Value = 19.63762386972934
$ws.Range("I5").Value = 27.83983774986796
$ws.Range("J5").Value = 8.992313170704639
$ws.Range("M5").Value = 19.84968926824474

$ws.Range("B6").Value = 15.34595793359739
$ws.Range("C6").Value = 12.11504620855905
$ws.Range("D6").Value = 14.29809479495895
$ws.Range("E6").Value = 15.19302168057127
$ws.Range("G6").Value = 50.91677166666221
$ws.Range("H6").Value = 19.63886988487305
$ws.Range("I6").Value = 27.84242735386619
$ws.Range("J6").Value = 8.993227173722209
$ws.Range("M6").Value = 19.84547026207437

$ws.Range("B7").Value = 15.46565760761456
$ws.Range("C7").Value = 12.22798006544363
$ws.Range("D7").Value = 14.29833501974834
$ws.Range("E7").Value = 15.18530741948609
$ws.Range("G7").Value = 50.98207434216363
$ws.Range("H7").Value = 19.63061871781331
$ws.Range("I7").Value = 27.82505628263027
$ws.Range("J7").Value = 8.986949544335006
$ws.Range("M7").Value = 19.87514156457271

$ws.Range("B8").Value = 15.99190527561428
$ws.Range("C8").Value = 12.72193718869632
$ws.Range("D8").Value = 14.30594832051716
$ws.Range("E8").Value = 15.1570828345307
$ws.Range("G8").Value = 51.30693997454164
$ws.Range("H8").Value = 19.60407321426915
$ws.Range("I8").Value = 27.76339967838134
$ws.Range("J8").Value = 8.960891494908843
$ws.Range("M8").Value = 20.01550156208109

$ws.Range("B9").Value = 17.00611073078496
$ws.Range("C9").Value = 13.66436035307222
$ws.Range("D9").Value = 14.34438502744164
$ws.Range("E9").Value = 15.1227027254936
$ws.Range("G9").Value = 52.0732315088961
$ws.Range("H9").Value = 19.5877097324712
$ws.Range("I9").Value = 27.69694632277466
$ws.Range("J9").Value = 8.915670593109988
$ws.Range("M9").Value = 20.32313791572958

$ws.Range("B10").Value = 17.73163414655891
$ws.Range("C10").Value = 14.33247776906617
$ws.Range("D10").Value = 14.38688156452355
$ws.Range("E10").Value = 15.11027568704965
$ws.Range("G10").Value = 52.71198069265418
$ws.Range("H10").Value = 19.59762769242768
$ws.Range("I10").Value = 27.68177704499963
$ws.Range("J10").Value = 8.886003511704414
$ws.Range("M10").Value = 20.56744336279937

$ws.Range("B11").Value = 18.05577191799222
$ws.Range("C11").Value = 14.62962939336344
$ws.Range("D11").Value = 14.40928512102555
$ws.Range("E11").Value = 15.10741635564731
$ws.Range("G11").Value = 53.01822937042301
$ws.Range("H11").Value = 19.60693225382735
$ws.Range("I11").Value = 27.68226780728396
$ws.Range("J11").Value = 8.8732737673228
$ws.Range("M11").Value = 20.6822171292475

$ws.Range("B12").Value = 18.17754059090876
$ws.Range("C12").Value = 14.74106738044051
$ws.Range("D12").Value = 14.41820751971193
$ws.Range("E12").Value = 15.10673566756008
$ws.Range("G12").Value = 53.13637720249786
$ws.Range("H12").Value = 19.61114622856456
$ws.Range("I12").Value = 27.68352214886532
$ws.Range("J12").Value = 8.868563058126618
$ws.Range("M12").Value = 20.72617207719081

$ws.Range("B13").Value = 18.15136079741628
$ws.Range("C13").Value = 14.71711722183729
$ws.Range("D13").Value = 14.41626646547487
$ws.Range("E13").Value = 15.10686437976542
$ws.Range("G13").Value = 53.11083635376564
$ws.Range("H13").Value = 19.6102079463693
$ws.Range("I13").Value = 27.68320440049353
$ws.Range("J13").Value = 8.869572717563972
$ws.Range("M13").Value = 20.71668417841185

$ws.Range("B14").Value = 18.0658101298433
$ws.Range("C14").Value = 14.63881984881874
$ws.Range("D14").Value = 14.41001039985589
$ws.Range("E14").Value = 15.10735229543496
$ws.Range("G14").Value = 53.02790626770368
$ws.Range("H14").Value = 19.60726509494849
$ws.Range("I14").Value = 27.68234955882744
$ws.Range("J14").Value = 8.872884016645678
$ws.Range("M14").Value = 20.68582363436808

$ws.Range("B15").Value = 18.01327730883016
$ws.Range("C15").Value = 14.59071576640206
$ws.Range("D15").Value = 14.40623540699909
$ws.Range("E15").Value = 15.10770352753938
$ws.Range("G15").Value = 52.97739049029763
$ws.Range("H15").Value = 19.60555247228048
$ws.Range("I15").Value = 27.68196524811779
$ws.Range("J15").Value = 8.874926567019884
$ws.Range("M15").Value = 20.6669838925828

$ws.Range("B16").Value = 17.7103211239675
$ws.Range("C16").Value = 14.31291196671225
$ws.Range("D16").Value = 14.38547899963085
$ws.Range("E16").Value = 15.11051880164045
$ws.Range("G16").Value = 52.69227570015612
$ws.Range("H16").Value = 19.59711618826715
$ws.Range("I16").Value = 27.68189420424154
$ws.Range("J16").Value = 8.886850811302292
$ws.Range("M16").Value = 20.56001324329402

$ws.Range("B17").Value = 17.52286160410779
$ws.Range("C17").Value = 14.14066881587027
$ws.Range("D17").Value = 14.37353032586931
$ws.Range("E17").Value = 15.11296172776307
$ws.Range("G17").Value = 52.52132997850752
$ws.Range("H17").Value = 19.59316948560757
$ws.Range("I17").Value = 27.68374797062503
$ws.Range("J17").Value = 8.894361860897297
$ws.Range("M17").Value = 20.49530055376109

$ws.Range("B18").Value = 17.4144907634188
$ws.Range("C18").Value = 14.04096707233236
$ws.Range("D18").Value = 14.36694709671513
$ws.Range("E18").Value = 15.11462978934613
$ws.Range("G18").Value = 52.42448718932985
$ws.Range("H18").Value = 19.59135052566745
$ws.Range("I18").Value = 27.6855096270043
$ws.Range("J18").Value = 8.898754135861543
$ws.Range("M18").Value = 20.45842383881272

$ws.Range("B19").Value = 17.37770796201868
$ws.Range("C19").Value = 14.00710482581151
$ws.Range("D19").Value = 14.36476790659546
$ws.Range("E19").Value = 15.11523971623476
$ws.Range("G19").Value = 52.39195443171333
$ws.Range("H19").Value = 19.59081208053916
$ws.Range("I19").Value = 27.68622534824886
$ws.Range("J19").Value = 8.900253682523775
$ws.Range("M19").Value = 20.44599806848594

$ws.Range("B20").Value = 17.54287476354633
$ws.Range("C20").Value = 14.15907065536839
$ws.Range("D20").Value = 14.37477235843374
$ws.Range("E20").Value = 15.11267445635365
$ws.Range("G20").Value = 52.53937473605086
$ws.Range("H20").Value = 19.59354292266588
$ws.Range("I20").Value = 27.68347861786922
$ws.Range("J20").Value = 8.893554835676404
$ws.Range("M20").Value = 20.50215391348909

$ws.Range("B21").Value = 18.09096585476471
$ws.Range("C21").Value = 14.66184802099201
$ws.Range("D21").Value = 14.41183607948024
$ws.Range("E21").Value = 15.10719806861549
$ws.Range("G21").Value = 53.05220636628273
$ws.Range("H21").Value = 19.60811073279679
$ws.Range("I21").Value = 27.68257160662379
$ws.Range("J21").Value = 8.871908431711869
$ws.Range("M21").Value = 20.69487500223463

$ws.Range("B22").Value = 18.44343941491429
$ws.Range("C22").Value = 14.98406272245844
$ws.Range("D22").Value = 14.43861446148548
$ws.Range("E22").Value = 15.10596259452075
$ws.Range("G22").Value = 53.40002811567508
$ws.Range("H22").Value = 19.62165685851626
$ws.Range("I22").Value = 27.68820893204659
$ws.Range("J22").Value = 8.858400901010823
$ws.Range("M22").Value = 20.82368599119073

$ws.Range("B23").Value = 18.25588135733798
$ws.Range("C23").Value = 14.81270855216882
$ws.Range("D23").Value = 14.42408966830902
$ws.Range("E23").Value = 15.10640747664222
$ws.Range("G23").Value = 53.21325759724762
$ws.Range("H23").Value = 19.61405843006127
$ws.Range("I23").Value = 27.68462846655141
$ws.Range("J23").Value = 8.865551719412194
$ws.Range("M23").Value = 20.75468605330621

$ws.Range("B24").Value = 17.53382866667927
$ws.Range("C24").Value = 14.15075328340192
$ws.Range("D24").Value = 14.37420994375544
$ws.Range("E24").Value = 15.11280351073885
$ws.Range("G24").Value = 52.53121222182331
$ws.Range("H24").Value = 19.59337268995564
$ws.Range("I24").Value = 27.68359822495862
$ws.Range("J24").Value = 8.893919461019188
$ws.Range("M24").Value = 20.49905448720271

$ws.Range("B25").Value = 16.73458715075612
$ws.Range("C25").Value = 13.4131189236539
$ws.Range("D25").Value = 14.33147603184542
$ws.Range("E25").Value = 15.12975299375662
$ws.Range("G25").Value = 51.85237966034417
$ws.Range("H25").Value = 19.58829467659553
$ws.Range("I25").Value = 27.70904412648193
$ws.Range("J25").Value = 8.927277502836501
$ws.Range("M25").Value = 20.23659201922778
